$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2924280.2
$ws.Range("J6").Value = 305
$ws.Range("L6").Value = 915
$ws.Range("N6").Value = -1139
$ws.Range("H12").Value = 241.8
$ws.Range("J12").Value = 214
$ws.Range("L12").Value = 214
$ws.Range("N12").Value = -554
$ws.Range("H29").Value = 1574.5
$ws.Range("I29").Value = 149
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 447
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = -166
$ws.Range("N29").Value = -9562
$ws.Range("H58").Value = 1489.5454
$ws.Range("I58").Value = 1348.5
$ws.Range("J58").Value = 2900
$ws.Range("K58").Value = 4045.5
$ws.Range("L58").Value = 8700
$ws.Range("M58").Value = -3895.5
$ws.Range("N58").Value = -9000
$ws.Range("H87").Value = 37684.617
$ws.Range("J87").Value = 39990.91
$ws.Range("L87").Value = 39990.91
$ws.Range("N87").Value = -42486.91
$ws.Range("H90").Value = 37684.617
$ws.Range("J90").Value = 39990.91
$ws.Range("L90").Value = 119972.73
$ws.Range("N90").Value = -132452.73

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20616
$ws.Range("H74").Value = 14113.625
$ws.Range("I74").Value = 14113.625
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 14113.625
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -13239.625
$ws.Range("H77").Value = 14113.625
$ws.Range("I77").Value = 14113.625
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 70568.125
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -66200.125
$ws.Range("H88").Value = 2237.5
$ws.Range("I88").Value = 2375
$ws.Range("J88").Value = 2100
$ws.Range("K88").Value = 2375
$ws.Range("L88").Value = 2100
$ws.Range("M88").Value = -1969
$ws.Range("N88").Value = -2912
$ws.Range("H91").Value = 2237.5
$ws.Range("I91").Value = 2375
$ws.Range("J91").Value = 2100
$ws.Range("K91").Value = 2375
$ws.Range("L91").Value = 2100
$ws.Range("M91").Value = -971
$ws.Range("N91").Value = -4908

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 116.666664
$ws.Range("I16").Value = 116.666664
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 116.666664
$ws.Range("L16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("M16").Value = 53.333336
$ws.Range("H29").Value = 800
$ws.Range("I29").Value = 800
$ws.Range("K29").Value = 800
$ws.Range("M29").Value = -511
$ws.Range("H86").Value = 4500
$ws.Range("J86").Value = 4500
$ws.Range("L86").Value = 4500
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 4500
$ws.Range("J89").Value = 4500
$ws.Range("L89").Value = 22500
$ws.Range("N89").Value = -33732

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H58").Value = 2916.2856
$ws.Range("I58").Value = 2569
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 2569
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -2366
$ws.Range("N58").Value = -5406
$ws.Range("H62").Value = 1113322.2
$ws.Range("I62").Value = 1252112.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 1252112.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1251488.5
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 1113322.2
$ws.Range("I65").Value = 1252112.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 6260562.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -6257442.5
$ws.Range("N65").Value = -21240
$ws.Range("H136").Value = 2916.2856
$ws.Range("I136").Value = 2569
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7707
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5157
$ws.Range("N136").Value = -20100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1291
$ws.Range("I55").Value = 187
$ws.Range("J55").Value = 3499
$ws.Range("K55").Value = 561
$ws.Range("L55").Value = 10497
$ws.Range("M55").Value = -384
$ws.Range("N55").Value = -10851

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 469
$ws.Range("I9").Value = 469
$ws.Range("K9").Value = 469
$ws.Range("M9").Value = -299
$ws.Range("H70").Value = 10069.8
$ws.Range("I70").Value = 10577.556
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 10577.556
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -10307.556
$ws.Range("N70").Value = -6040
$ws.Range("H73").Value = 10069.8
$ws.Range("I73").Value = 10577.556
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 10577.556
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -9641.556
$ws.Range("N73").Value = -7372
$ws.Range("H113").Value = 12692.2
$ws.Range("I113").Value = 2365.25
$ws.Range("K113").Value = 2365.25
$ws.Range("M113").Value = -195.25
$ws.Range("H126").Value = 2689.5833
$ws.Range("I126").Value = 2427.56
$ws.Range("J126").Value = 3285.0908
$ws.Range("K126").Value = 7282.68
$ws.Range("L126").Value = 9855.2724
$ws.Range("M126").Value = -4812.68
$ws.Range("N126").Value = -14795.2724

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 833983.25
$ws.Range("I22").Value = 1111455.5
$ws.Range("J22").Value = 1566.6666
$ws.Range("K22").Value = 1111455.5
$ws.Range("L22").Value = 1566.6666
$ws.Range("M22").Value = -1111160.5
$ws.Range("N22").Value = -2156.6666
$ws.Range("H27").Value = 833983.25
$ws.Range("I27").Value = 1111455.5
$ws.Range("J27").Value = 1566.6666
$ws.Range("K27").Value = 1111455.5
$ws.Range("L27").Value = 1566.6666
$ws.Range("M27").Value = -1111348.5
$ws.Range("N27").Value = -1780.6666
$ws.Range("H132").Value = 4609.364
$ws.Range("I132").Value = 5013.3335
$ws.Range("J132").Value = 3743.7144
$ws.Range("K132").Value = 15040.0005
$ws.Range("L132").Value = 11231.1432
$ws.Range("M132").Value = -12510.0005
$ws.Range("N132").Value = -16291.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 447.625
$ws.Range("I113").Value = 420
$ws.Range("J113").Value = 493.66666
$ws.Range("K113").Value = 1260
$ws.Range("L113").Value = 1480.99998
$ws.Range("M113").Value = 910
$ws.Range("N113").Value = -5820.999980000001
$ws.Range("H132").Value = 3214.8572
$ws.Range("I132").Value = 3376
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 10128
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -7598
$ws.Range("N132").Value = -14060
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120
$ws.Range("H139").Value = 33636.25
$ws.Range("I139").Value = 59650
$ws.Range("J139").Value = 29920
$ws.Range("K139").Value = 59650
$ws.Range("L139").Value = 29920
$ws.Range("M139").Value = -54510
$ws.Range("N139").Value = -40200
